$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: reformat earthquake timestamps from "Www Mon dd HH:MM:SS CST yyyy"
#     to "dd-mm-yyyy HH:MM:SS" (kept as plain text, same as before the edit) ---
$ws.Range("A2").Value = "20-08-2020 22:58:39"
$ws.Range("A3").Value = "19-08-2017 21:57:38"
$ws.Range("A4").Value = "09-09-2007 11:47:28"
$ws.Range("A5").Value = "19-01-2018 22:57:38"
$ws.Range("A6").Value = "14-08-2012 16:52:33"
$ws.Range("A7").Value = "05-04-2003 07:43:24"
$ws.Range("A8").Value = "05-03-2018 07:43:24"

# The date cells used a date number format (numFmtId 14); now that the
# column holds plain "dd-mm-yyyy HH:MM:SS" text, switch it to the same
# plain-text format already used elsewhere in the sheet (column F/G),
# so the now-unused date format is dropped from the workbook.
$txtFormat = $ws.Range("G6").NumberFormat
$ws.Range("A2:A8").NumberFormat = $txtFormat

# --- Column H: shorten/clean up the "Lugar" detail text ---
$ws.Range("H2").Value = "50 km al Suroeste de Jacó en Playa Hermosa de Garabito"
$ws.Range("H3").Value = "6.13 km hacia el Noreste de Puerto Armuelles"
$ws.Range("H4").Value = "3.68 km hacia el Noreste de Dos Rios de Osa de Puntarenas"
$ws.Range("H5").Value = "3 km al SE de Paso Real en Buenos Aires además Zona sur y  Valle Central"
$ws.Range("H6").Value = "9.85 km al SO de las Juntas de Abangares en Nicoya y Santa Cruz."
$ws.Range("H7").Value = "Orosí 3.9 Km SE en El Carmen y Tres Ríos"
$ws.Range("H8").Value = "Punta Banco en  Golfito 2.6 Km de Bijagua de Upala en Fuerte en Punta Banco de Golfito"

# --- Update the last saved selection to match the author's final click ---
$ws.Range("A8").Select()
